$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the B2:D11 block as text first so the numeric-looking values
# are stored as literal strings (matching the source table), not numbers.
$ws.Range("B2:D11").NumberFormat = "@"

$ws.Range("B2").Value = "3.0"
$ws.Range("C2").Value = "-17.6902860206768"
$ws.Range("D2").Value = "1.000005"
$ws.Range("B3").Value = "20.6902860206768"
$ws.Range("C3").Value = "-13.5173641896727"
$ws.Range("D3").Value = "0.855004421060107"
$ws.Range("B4").Value = "34.2076502103495"
$ws.Range("C4").Value = "-3.51949495555342"
$ws.Range("D4").Value = "0.395156174322171"
$ws.Range("B5").Value = "37.7271451659029"
$ws.Range("C5").Value = "-0.685513898070205"
$ws.Range("D5").Value = "0.0932881335196885"
$ws.Range("B6").Value = "38.4126590639731"
$ws.Range("C6").Value = "-0.126050414764144"
$ws.Range("D6").Value = "0.0178460412471977"
$ws.Range("B7").Value = "38.5387094787372"
$ws.Range("C7").Value = "-0.0229327615412593"
$ws.Range("D7").Value = "0.0032707482027567"
$ws.Range("B8").Value = "38.5616422402785"
$ws.Range("C8").Value = "-0.0041641661015425"
$ws.Range("D8").Value = "0.0005947039651052"
$ws.Range("B9").Value = "38.56580640638"
$ws.Range("C9").Value = "-0.0007558700112753"
$ws.Range("D9").Value = "0.0001079756003975"
$ws.Range("B10").Value = "38.5665622763913"
$ws.Range("C10").Value = "-0.000137195062706"
$ws.Range("D10").Value = "1.9599102607546e-05"
$ws.Range("B11").Value = "38.566699471454"
$ws.Range("C11").Value = "-2.49014605948616e-05"
$ws.Range("D11").Value = "3.55734518603515e-06"

# Drop the temporary text format again so the cells keep the workbook
# default style (no explicit number format), matching the source file.
$ws.Range("B2:D11").ClearFormats()

# The root-finding table now only needs 10 iterations (rows 2-11) to
# converge, so drop the old extra rows 12-16.
$ws.Rows("12:16").Delete()
